# Auto-generated edit script applying the Tonberry_Profits.xlsx numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1255.4286
$ws.Range("J88").Value = 947.5
$ws.Range("L88").Value = 947.5
$ws.Range("N88").Value = -1759.5
$ws.Range("H91").Value = 1255.4286
$ws.Range("J91").Value = 947.5
$ws.Range("L91").Value = 947.5
$ws.Range("N91").Value = -3755.5
$ws.Range("H137").Value = 35709.414
$ws.Range("I137").Value = 805.9167
$ws.Range("K137").Value = 2417.7501
$ws.Range("M137").Value = 132.2498999999998
$ws.Range("H138").Value = 3173.508
$ws.Range("I138").Value = 4641.7144
$ws.Range("J138").Value = 2754.0205
$ws.Range("K138").Value = 13925.1432
$ws.Range("L138").Value = 8262.0615
$ws.Range("M138").Value = -8785.143199999999
$ws.Range("N138").Value = -18542.0615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2487.204
$ws.Range("I32").Value = 1963.7701
$ws.Range("K32").Value = 1963.7701
$ws.Range("M32").Value = -1676.7701
$ws.Range("H61").Value = 2185.7407
$ws.Range("I61").Value = 692.8125
$ws.Range("K61").Value = 692.8125
$ws.Range("M61").Value = -480.8125
$ws.Range("H88").Value = 3847.5
$ws.Range("I88").Value = 1980
$ws.Range("J88").Value = 4114.2856
$ws.Range("K88").Value = 1980
$ws.Range("L88").Value = 4114.2856
$ws.Range("M88").Value = -1574
$ws.Range("N88").Value = -4926.2856
$ws.Range("H91").Value = 3847.5
$ws.Range("I91").Value = 1980
$ws.Range("J91").Value = 4114.2856
$ws.Range("K91").Value = 1980
$ws.Range("L91").Value = 4114.2856
$ws.Range("M91").Value = -576
$ws.Range("N91").Value = -6922.2856
$ws.Range("H97").Value = 2313
$ws.Range("I97").Value = 2313
$ws.Range("K97").Value = 2313
$ws.Range("M97").Value = -1817
$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("M102").Value = 122
$ws.Range("H135").Value = 49890
$ws.Range("J135").Value = 49890
$ws.Range("L135").Value = 49890
$ws.Range("N135").Value = -60030
$ws.Range("H136").Value = 2185.7407
$ws.Range("I136").Value = 692.8125
$ws.Range("K136").Value = 2078.4375
$ws.Range("M136").Value = 471.5625
$ws.Range("H139").Value = 51197
$ws.Range("J139").Value = 51197
$ws.Range("L139").Value = 51197
$ws.Range("N139").Value = -61477

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2188.3076
$ws.Range("I20").Value = 1853.8
$ws.Range("J20").Value = 3303.3333
$ws.Range("K20").Value = 1853.8
$ws.Range("L20").Value = 3303.3333
$ws.Range("M20").Value = -1606.8
$ws.Range("N20").Value = -3797.3333
$ws.Range("H94").Value = 736.44446
$ws.Range("I94").Value = 673.2308
$ws.Range("J94").Value = 900.8
$ws.Range("K94").Value = 673.2308
$ws.Range("L94").Value = 900.8
$ws.Range("M94").Value = -222.2308
$ws.Range("N94").Value = -1802.8
$ws.Range("H107").Value = 1907.6923
$ws.Range("I107").Value = 1528.1818
$ws.Range("J107").Value = 3995
$ws.Range("K107").Value = 1528.1818
$ws.Range("L107").Value = 3995
$ws.Range("M107").Value = 391.8181999999999
$ws.Range("N107").Value = -7835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2583.3333
$ws.Range("J31").Value = 3500.6667
$ws.Range("L31").Value = 3500.6667
$ws.Range("N31").Value = -4090.6667
$ws.Range("H34").Value = 2583.3333
$ws.Range("J34").Value = 3500.6667
$ws.Range("L34").Value = 3500.6667
$ws.Range("N34").Value = -3904.6667
$ws.Range("H107").Value = 386.4
$ws.Range("I107").Value = 349.66666
$ws.Range("K107").Value = 349.66666
$ws.Range("M107").Value = 1570.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 292.60605
$ws.Range("I5").Value = 229.2
$ws.Range("J5").Value = 345.44446
$ws.Range("K5").Value = 687.5999999999999
$ws.Range("L5").Value = 1036.33338
$ws.Range("M5").Value = -575.5999999999999
$ws.Range("N5").Value = -1260.33338
$ws.Range("H68").Value = 2448.762
$ws.Range("J68").Value = 2806.1875
$ws.Range("L68").Value = 8418.5625
$ws.Range("N68").Value = -10040.5625
$ws.Range("H71").Value = 2448.762
$ws.Range("J71").Value = 2806.1875
$ws.Range("L71").Value = 25255.6875
$ws.Range("N71").Value = -33367.6875
$ws.Range("H92").Value = 422.22223
$ws.Range("J92").Value = 422.22223
$ws.Range("L92").Value = 1266.66669
$ws.Range("N92").Value = -3762.66669
$ws.Range("H107").Value = 2438.889
$ws.Range("J107").Value = 2510.3684
$ws.Range("L107").Value = 7531.1052
$ws.Range("N107").Value = -11371.1052
$ws.Range("H113").Value = 68023.336
$ws.Range("J113").Value = 818.63635
$ws.Range("L113").Value = 2455.90905
$ws.Range("N113").Value = -6795.90905
$ws.Range("H132").Value = 2166
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2166
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 19494
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -24554
$ws.Range("H135").Value = 292.60605
$ws.Range("I135").Value = 229.2
$ws.Range("J135").Value = 345.44446
$ws.Range("K135").Value = 2062.8
$ws.Range("L135").Value = 3109.00014
$ws.Range("M135").Value = 472.2000000000003
$ws.Range("N135").Value = -8179.00014
$ws.Range("H137").Value = 3466.1516
$ws.Range("J137").Value = 5143
$ws.Range("L137").Value = 15429
$ws.Range("N137").Value = -25629

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 758.3333
$ws.Range("I107").Value = 201
$ws.Range("J107").Value = 1204.2
$ws.Range("K107").Value = 201
$ws.Range("L107").Value = 1204.2
$ws.Range("M107").Value = 1719
$ws.Range("N107").Value = -5044.2
$ws.Range("H113").Value = 1733.3334
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -5940
$ws.Range("H126").Value = 1349120.2
$ws.Range("I126").Value = 1854530.5
$ws.Range("J126").Value = 85594.5
$ws.Range("K126").Value = 5563591.5
$ws.Range("L126").Value = 256783.5
$ws.Range("M126").Value = -5561121.5
$ws.Range("N126").Value = -261723.5
$ws.Range("H132").Value = 1482762.6
$ws.Range("I132").Value = 12823229
$ws.Range("J132").Value = 3571.348
$ws.Range("K132").Value = 38469687
$ws.Range("L132").Value = 10714.044
$ws.Range("M132").Value = -38467157
$ws.Range("N132").Value = -15774.044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2167.8333
$ws.Range("J22").Value = 2167.8333
$ws.Range("L22").Value = 2167.8333
$ws.Range("N22").Value = -2757.8333
$ws.Range("H27").Value = 2167.8333
$ws.Range("J27").Value = 2167.8333
$ws.Range("L27").Value = 2167.8333
$ws.Range("N27").Value = -2381.8333
$ws.Range("H40").Value = 4943.364
$ws.Range("I40").Value = 4769.5
$ws.Range("J40").Value = 5042.7144
$ws.Range("K40").Value = 4769.5
$ws.Range("L40").Value = 5042.7144
$ws.Range("M40").Value = -4633.5
$ws.Range("N40").Value = -5314.7144
$ws.Range("H122").Value = 5130.1904
$ws.Range("I122").Value = 4210.385
$ws.Range("K122").Value = 12631.155
$ws.Range("M122").Value = -10181.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 266713380
$ws.Range("J47").Value = 266713380
$ws.Range("L47").Value = 266713380
$ws.Range("N47").Value = -266714524
$ws.Range("H81").Value = 3149.7778
$ws.Range("I81").Value = 2630.4614
$ws.Range("K81").Value = 5260.9228
$ws.Range("M81").Value = -4199.9228
$ws.Range("H84").Value = 3149.7778
$ws.Range("I84").Value = 2630.4614
$ws.Range("K84").Value = 26304.614
$ws.Range("M84").Value = -21000.614
$ws.Range("H141").Value = 72442.836
$ws.Range("J141").Value = 72442.836
$ws.Range("L141").Value = 72442.836
$ws.Range("N141").Value = -82802.836

Write-Output "Applied 202 cell updates across 8 sheets"
